$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = [double]"1.724903831991619E-08"
$ws.Range("E2").Value = [double]"1.724903831991619E-08"

# Row 3
$ws.Range("D3").Value = [double]"3.209255032175693E-05"
$ws.Range("E3").Value = [double]"3.209255032175693E-05"

# Row 4
$ws.Range("D4").Value = 0.9999999946810769
$ws.Range("E4").Value = 0.9999999946810769

# Row 5
$ws.Range("D5").Value = 0.9999935213091597
$ws.Range("E5").Value = 0.9999935213091597

# Row 6
$ws.Range("D6").Value = 0.9999999908814255
$ws.Range("E6").Value = 0.9999999908814255

# Row 7
$ws.Range("D7").Value = [double]"3.37352351957873E-07"
$ws.Range("E7").Value = 0.9999996626476481

# Row 8
$ws.Range("D8").Value = 0.9999856222915967
$ws.Range("E8").Value = [double]"1.437770840329033E-05"

# Row 9
$ws.Range("C9").Value = $true
$ws.Range("D9").Value = 0.8397839618747578
$ws.Range("E9").Value = 0.1602160381252422

# Row 10
$ws.Range("D10").Value = [double]"2.180930155392142E-09"
$ws.Range("E10").Value = 0.9999999978190699

# Row 11
$ws.Range("D11").Value = [double]"1.646093504184164E-05"
$ws.Range("E11").Value = 0.9999835390649582
$ws.Range("F11").Value = 9.554677963256836
$ws.Range("G11").Value = 0.4
